$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 93
$ws1.Range("F4").Value = 1409
$ws1.Range("F5").Value = 166
$ws1.Range("F6").Value = 31
$ws1.Range("F8").Value = 9551
$ws1.Range("F9").Value = 158
$ws1.Range("F10").Value = 106
$ws1.Range("F11").Value = 231
$ws1.Range("F12").Value = 184
$ws1.Range("F14").Value = 6550
$ws1.Range("F15").Value = 1075
$ws1.Range("F16").Value = 109
$ws1.Range("F17").Value = 47
$ws1.Range("F18").Value = 165

# --- Sheet "演出" (Performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 41

# --- Sheet "全部类型" (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 93
$ws4.Range("F4").Value = 1409
$ws4.Range("F5").Value = 166
$ws4.Range("F6").Value = 31
$ws4.Range("F8").Value = 41
$ws4.Range("F10").Value = 9551
$ws4.Range("F11").Value = 158
$ws4.Range("F12").Value = 106
$ws4.Range("F13").Value = 231
$ws4.Range("F14").Value = 184
$ws4.Range("F16").Value = 6550
$ws4.Range("F17").Value = 1075
$ws4.Range("F18").Value = 109
$ws4.Range("F19").Value = 47
$ws4.Range("F20").Value = 165
